# Module 6 - Volt-VAR Optimization: fix "ANCI" -> "ANSI" typo (two spots) and
# tidy up the run split on the capacitor-switching-thresholds bullet (slide 19).

$p = $ppt.ActivePresentation

# --- Slide 3: body text "...limits set by ANCI C84.1, 120V +/-5%." ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)          # "Content Placeholder 2"
$tr3 = $sh3.TextFrame.TextRange
$word3 = $tr3.Characters(87, 5)    # "ANCI " (including trailing space)
$word3.Text = "ANSI "

# --- Slide 7: title "ANCI C84.1" ---
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(1)          # "Title 1"
$tr7 = $sh7.TextFrame.TextRange
$word7 = $tr7.Characters(1, 5)     # "ANCI " (including trailing space)
$word7.Text = "ANSI "

# --- Slide 19: merge the trailing "<tab>ratio constants <en dash> create
#     hysteresis-like condition" runs into a single run ---
$s19 = $p.Slides.Item(19)
$sh19 = $s19.Shapes.Item(2)        # "Content Placeholder 2"
$tr19 = $sh19.TextFrame.TextRange
$tail19 = $tr19.Characters(356, 51)
$tail19.Text = "`tratio constants " + [char]0x2013 + " create hysteresis-like condition"
